# Auto-applied data refresh for Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per the
# latest scheduled market-data pull. Some profit cells are cleared when the
# recomputed value is no longer negative/applicable; a couple of previously
# blank profit cells now receive a value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 533.0526
$ws.Range("I33").Value = 202.83333
$ws.Range("J33").Value = 1099.1428
$ws.Range("K33").Value = 202.83333
$ws.Range("L33").Value = 1099.1428
$ws.Range("M33").Value = 26.16667000000001
$ws.Range("N33").Value = -1557.1428
$ws.Range("H63").Value = 100000
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101248
$ws.Range("H66").Value = 100000
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306240
$ws.Range("H106").Value = 6317.375
$ws.Range("I106").Value = 5941.2856
$ws.Range("K106").Value = 5941.2856
$ws.Range("M106").Value = -5310.2856
$ws.Range("H132").Value = 17800.766
$ws.Range("I132").Value = 2341.3914
$ws.Range("K132").Value = 7024.174199999999
$ws.Range("M132").Value = -4494.174199999999
$ws.Range("H138").Value = 3573.6382
$ws.Range("J138").Value = 3766.4358
$ws.Range("L138").Value = 11299.3074
$ws.Range("N138").Value = -21579.3074

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H61").Value = 1282.2142
$ws.Range("I61").Value = 1337.25
$ws.Range("K61").Value = 1337.25
$ws.Range("M61").Value = -1125.25
$ws.Range("H110").Value = 5225.8667
$ws.Range("I110").Value = 5463
$ws.Range("J110").Value = 4573.75
$ws.Range("K110").Value = 5463
$ws.Range("L110").Value = 4573.75
$ws.Range("M110").Value = -3418
$ws.Range("N110").Value = -8663.75
$ws.Range("H122").Value = 3821.4348
$ws.Range("I122").Value = 3240.5356
$ws.Range("J122").Value = 4725.0557
$ws.Range("K122").Value = 9721.606800000001
$ws.Range("L122").Value = 14175.1671
$ws.Range("M122").Value = -7271.606800000001
$ws.Range("N122").Value = -19075.1671
$ws.Range("H124").Value = 38321.668
$ws.Range("J124").Value = 38321.668
$ws.Range("L124").Value = 38321.668
$ws.Range("N124").Value = -48141.668
$ws.Range("H125").Value = 69909.09
$ws.Range("J125").Value = 69909.09
$ws.Range("L125").Value = 69909.09
$ws.Range("N125").Value = -79749.09
$ws.Range("H136").Value = 1282.2142
$ws.Range("I136").Value = 1337.25
$ws.Range("K136").Value = 4011.75
$ws.Range("M136").Value = -1461.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1622.6875
$ws.Range("I20").Value = 1862.7778
$ws.Range("K20").Value = 1862.7778
$ws.Range("M20").Value = -1615.7778
$ws.Range("H80").Value = 172.25
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 172.25
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H94").Value = 2548
$ws.Range("I94").Value = 2277.5
$ws.Range("J94").Value = 3209.2222
$ws.Range("K94").Value = 2277.5
$ws.Range("L94").Value = 3209.2222
$ws.Range("M94").Value = -1826.5
$ws.Range("N94").Value = -4111.2222
$ws.Range("H105").Value = 2091.25
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H116").Value = 69997.5
$ws.Range("J116").Value = 69997.5
$ws.Range("L116").Value = 69997.5
$ws.Range("N116").Value = -79175.5
$ws.Range("H134").Value = 1971.5
$ws.Range("I134").Value = 1681.0952
$ws.Range("J134").Value = 4004.3333
$ws.Range("K134").Value = 5043.2856
$ws.Range("L134").Value = 12012.9999
$ws.Range("M134").Value = -2508.2856
$ws.Range("N134").Value = -17082.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1750
$ws.Range("I12").Value = 1750
$ws.Range("K12").Value = 1750
$ws.Range("M12").Value = -1580
$ws.Range("H58").Value = 1315.2646
$ws.Range("I58").Value = 668.5357
$ws.Range("J58").Value = 4333.3335
$ws.Range("K58").Value = 668.5357
$ws.Range("L58").Value = 4333.3335
$ws.Range("M58").Value = -465.5357
$ws.Range("N58").Value = -4739.3335
$ws.Range("H104").Value = 54996.5
$ws.Range("J104").Value = 54996.5
$ws.Range("L104").Value = 54996.5
$ws.Range("N104").Value = -60238.5
$ws.Range("H136").Value = 1315.2646
$ws.Range("I136").Value = 668.5357
$ws.Range("J136").Value = 4333.3335
$ws.Range("K136").Value = 2005.6071
$ws.Range("L136").Value = 13000.0005
$ws.Range("M136").Value = 544.3928999999998
$ws.Range("N136").Value = -18100.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 334417.34
$ws.Range("I51").Value = 334417.34
$ws.Range("K51").Value = 1003252.02
$ws.Range("M51").Value = -1002792.02
$ws.Range("H125").Value = 2100
$ws.Range("I125").Value = 2900
$ws.Range("K125").Value = 8700
$ws.Range("M125").Value = -3780
$ws.Range("H137").Value = 1650.1818
$ws.Range("I137").Value = 1487.7
$ws.Range("J137").Value = 3275
$ws.Range("K137").Value = 4463.1
$ws.Range("L137").Value = 9825
$ws.Range("M137").Value = 636.8999999999996
$ws.Range("N137").Value = -20025

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 192456.17
$ws.Range("I70").Value = 373239.66
$ws.Range("J70").Value = 11672.667
$ws.Range("K70").Value = 373239.66
$ws.Range("L70").Value = 11672.667
$ws.Range("M70").Value = -372969.66
$ws.Range("N70").Value = -12212.667
$ws.Range("H73").Value = 192456.17
$ws.Range("I73").Value = 373239.66
$ws.Range("J73").Value = 11672.667
$ws.Range("K73").Value = 373239.66
$ws.Range("L73").Value = 11672.667
$ws.Range("M73").Value = -372303.66
$ws.Range("N73").Value = -13544.667
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 710660
$ws.Range("I43").Value = 9900
$ws.Range("J43").Value = 760714.3
$ws.Range("K43").Value = 9900
$ws.Range("L43").Value = 760714.3
$ws.Range("M43").Value = -9707
$ws.Range("N43").Value = -761100.3
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 1900.2142
$ws.Range("I61").Value = 1865.4
$ws.Range("K61").Value = 1865.4
$ws.Range("M61").Value = -1663.4
$ws.Range("H100").Value = 190680.67
$ws.Range("I100").Value = 374361.34
$ws.Range("K100").Value = 374361.34
$ws.Range("M100").Value = -373820.34
$ws.Range("H113").Value = 1900.2142
$ws.Range("I113").Value = 1865.4
$ws.Range("K113").Value = 1865.4
$ws.Range("M113").Value = 304.5999999999999
$ws.Range("H120").Value = 50349
$ws.Range("J120").Value = 50349
$ws.Range("L120").Value = 50349
$ws.Range("N120").Value = -60025
$ws.Range("H122").Value = 5263.227
$ws.Range("I122").Value = 4808.0625
$ws.Range("J122").Value = 6477
$ws.Range("K122").Value = 14424.1875
$ws.Range("L122").Value = 19431
$ws.Range("M122").Value = -11974.1875
$ws.Range("N122").Value = -24331
$ws.Range("H127").Value = 69576.92
$ws.Range("J127").Value = 69576.92
$ws.Range("L127").Value = 69576.92
$ws.Range("N127").Value = -79496.92
$ws.Range("H132").Value = 2506.4
$ws.Range("I132").Value = 2112
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 6336
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -3806
$ws.Range("N132").Value = -35060
$ws.Range("H137").Value = 41700
$ws.Range("J137").Value = 39666.668
$ws.Range("L137").Value = 39666.668
$ws.Range("N137").Value = -49866.668
$ws.Range("H139").Value = 38054
$ws.Range("J139").Value = 38054
$ws.Range("L139").Value = 38054
$ws.Range("N139").Value = -48334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 29498.334
$ws.Range("J54").Value = 38495
$ws.Range("L54").Value = 38495
$ws.Range("N54").Value = -39535
